$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $cell.Value = "'" + $val
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "50.994.26"
$ws.Range("E2").Value = "  -1.56%  "
$ws.Range("D3").Value = "2.939.80"
$ws.Range("E3").Value = "  -2.14%  "
$ws.Range("E4").Value = "  +0.28%  "
Set-TextValue $ws.Range("D5") "375.65"
$ws.Range("E5").Value = "  -2.02%  "
Set-TextValue $ws.Range("D6") "100.95"
$ws.Range("E6").Value = "  -3.90%  "
$ws.Range("E7").Value = "  -1.55%  "
Set-TextValue $ws.Range("D9") "0.588"
$ws.Range("E9").Value = "  -1.82%  "
Set-TextValue $ws.Range("D10") "36.29"
$ws.Range("E10").Value = "  -3.39%  "
$ws.Range("E11").Value = "  -0.89%  "
Set-TextValue $ws.Range("D12") "0.0850"
$ws.Range("E12").Value = "  -0.02%  "
$ws.Range("D13").Value = "3.408.61"
$ws.Range("E13").Value = "  -1.39%  "
Set-TextValue $ws.Range("D14") "18.10"
$ws.Range("E14").Value = "  -2.23%  "
Set-TextValue $ws.Range("D15") "7.57"
$ws.Range("E15").Value = "  -0.11%  "
$ws.Range("D16").Value = "2.941.20"
$ws.Range("E16").Value = "  -1.77%  "
Set-TextValue $ws.Range("D17") "0.995"
$ws.Range("E17").Value = "  +0.62%  "
$ws.Range("E18").Value = "  +47.01%  "
$ws.Range("D19").Value = "50.963.33"
$ws.Range("E19").Value = "  -1.28%  "
$ws.Range("E20").Value = "  -7.33%  "
Set-TextValue $ws.Range("D21") "12.45"
$ws.Range("E21").Value = "  -4.41%  "
$ws.Range("D22").Value = "0.0₃0954"
$ws.Range("E22").Value = "  -1.41%  "
Set-TextValue $ws.Range("D23") "265.09"
$ws.Range("E23").Value = "  +0.47%  "
Set-TextValue $ws.Range("D24") "68.72"
$ws.Range("E24").Value = "  -0.63%  "
Set-TextValue $ws.Range("D25") "3.13"
$ws.Range("E25").Value = "  +6.70%  "
Set-TextValue $ws.Range("D26") "8.11"
$ws.Range("E26").Value = "  -2.89%  "
Set-TextValue $ws.Range("D27") "7.54"
$ws.Range("E27").Value = "  -2.58%  "
Set-TextValue $ws.Range("D28") "0.999"
$ws.Range("E28").Value = "  -0.12%  "
Set-TextValue $ws.Range("D31") "0.109"
$ws.Range("E31").Value = "  -5.97%  "
Set-TextValue $ws.Range("D32") "9.99"
$ws.Range("E32").Value = "  +0.68%  "
Set-TextValue $ws.Range("D33") "50.69"
$ws.Range("E33").Value = "  -0.74%  "
$ws.Range("E34").Value = "  -0.88%  "
Set-TextValue $ws.Range("D35") "33.39"
$ws.Range("E35").Value = "  -4.94%  "
Set-TextValue $ws.Range("D36") "0.0441"
$ws.Range("E36").Value = "  -3.47%  "
$ws.Range("E37").Value = "  +0.17%  "
Set-TextValue $ws.Range("D38") "3.15"
$ws.Range("E38").Value = "  +3.27%  "
$ws.Range("E39").Value = "  -0.62%  "
$ws.Range("E40").Value = "  -4.43%  "
Set-TextValue $ws.Range("D41") "1.79"
$ws.Range("E41").Value = "  -3.72%  "
Set-TextValue $ws.Range("D42") "2.47"
$ws.Range("E42").Value = "  -5.29%  "
Set-TextValue $ws.Range("D43") "120.32"
$ws.Range("E43").Value = "  -1.66%  "
$ws.Range("E44").Value = "  -2.32%  "
Set-TextValue $ws.Range("D45") "3.36"
$ws.Range("E45").Value = "  +1.60%  "
$ws.Range("E46").Value = "  -0.57%  "
$ws.Range("E47").Value = "  -3.76%  "
$ws.Range("E48").Value = "  -2.15%  "
$ws.Range("D49").Value = "1.986.91"
$ws.Range("E49").Value = "  -2.72%  "
Set-TextValue $ws.Range("D50") "0.0330"
$ws.Range("E50").Value = "  -1.44%  "
Set-TextValue $ws.Range("D51") "1.31"
$ws.Range("E51").Value = "  +1.35%  "

# Row 29/30 swap with updated values
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D29") "25.61"
$ws.Range("E29").Value = "  -2.10%  "
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D30") "0.163"
$ws.Range("E30").Value = "  -4.77%  "
